# Update cryptocurrency price snapshot (column D) with the latest figures.
# Values in column D are stored as plain text (e.g. "245.90", "0.03123")
# rather than numbers, so each one is written back as text - forcing a
# Text number format for the write and then restoring the cell's original
# ("Normal") style - to keep exact textual formatting such as trailing
# zeros and very small decimals instead of having Excel coerce the text
# into a number (which would normalize/round it and could flip it to
# scientific notation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param(
        [string]$Address,
        [string]$Text
    )

    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

Set-PriceText "D2"  "245.87"
Set-PriceText "D3"  "22.00"
Set-PriceText "D4"  "5.444"
Set-PriceText "D6"  "3.424"
Set-PriceText "D7"  "6.349"
Set-PriceText "D8"  "0.8177"
Set-PriceText "D9"  "1.040"
Set-PriceText "D12" "0.03117"
Set-PriceText "D13" "0.03117"
Set-PriceText "D15" "0.09377"
Set-PriceText "D16" "0.001609"
Set-PriceText "D17" "0.04816"
Set-PriceText "D18" "0.0005849"
Set-PriceText "D19" "0.006291"
Set-PriceText "D20" "0.004130"
Set-PriceText "D21" "0.0009935"
Set-PriceText "D23" "3.741"
Set-PriceText "D24" "2.198"
Set-PriceText "D25" "0.3243"
Set-PriceText "D26" "0.1330"
Set-PriceText "D27" "0.0003998"
Set-PriceText "D40" "0.03882"
Set-PriceText "D41" "0.006669"
Set-PriceText "D42" "0.1072"
Set-PriceText "D43" "0.002639"
Set-PriceText "D44" "0.006595"
Set-PriceText "D45" "0.00005614"
Set-PriceText "D46" "0.00000000750"
Set-PriceText "D47" "0.3899"
Set-PriceText "D49" "0.00002099"
Set-PriceText "D50" "0.01010"
